$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1024.25
$ws.Range("I18").Value = 366.5
$ws.Range("K18").Value = 366.5
$ws.Range("M18").Value = -82.5

$ws.Range("H28").Value = 57175
$ws.Range("I28").Value = 112473.336
$ws.Range("K28").Value = 112473.336
$ws.Range("M28").Value = -111988.336

$ws.Range("H31").Value = 999
$ws.Range("J31").Value = 999
$ws.Range("L31").Value = 2997
$ws.Range("N31").Value = -3457

$ws.Range("H70").Value = 101953.6
$ws.Range("I70").Value = 2064.75
$ws.Range("J70").Value = 168546.17
$ws.Range("K70").Value = 6194.25
$ws.Range("L70").Value = 505638.51
$ws.Range("M70").Value = -5924.25
$ws.Range("N70").Value = -506178.51

$ws.Range("H73").Value = 101953.6
$ws.Range("I73").Value = 2064.75
$ws.Range("J73").Value = 168546.17
$ws.Range("K73").Value = 6194.25
$ws.Range("L73").Value = 505638.51
$ws.Range("M73").Value = -5258.25
$ws.Range("N73").Value = -507510.51

$ws.Range("H76").Value = 78013.57000000001
$ws.Range("I76").Value = 90181.914
$ws.Range("K76").Value = 90181.914
$ws.Range("M76").Value = -89866.914

$ws.Range("H79").Value = 78013.57000000001
$ws.Range("I79").Value = 90181.914
$ws.Range("K79").Value = 90181.914
$ws.Range("M79").Value = -89089.914

$ws.Range("H87").Value = 74957.14
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 74957.14
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H100").Value = 7457.875
$ws.Range("I100").Value = 2703.7778
$ws.Range("K100").Value = 2703.7778
$ws.Range("M100").Value = -2162.7778

$ws.Range("H125").Value = 8041.9
$ws.Range("I125").Value = 7546.846
$ws.Range("J125").Value = 8961.286
$ws.Range("K125").Value = 67921.614
$ws.Range("L125").Value = 80651.57399999999
$ws.Range("M125").Value = -65461.614
$ws.Range("N125").Value = -85571.57399999999

$ws.Range("H127").Value = 9652.214
$ws.Range("I127").Value = 11303.091
$ws.Range("K127").Value = 33909.273
$ws.Range("M127").Value = -28949.273

$ws.Range("H131").Value = 2918.9473
$ws.Range("I131").Value = 2576.8572
$ws.Range("K131").Value = 7730.571599999999
$ws.Range("M131").Value = -2690.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 94
$ws.Range("I5").Value = 92.28570999999999
$ws.Range("K5").Value = 92.28570999999999
$ws.Range("M5").Value = 19.71429000000001

$ws.Range("H32").Value = 3775.6726
$ws.Range("I32").Value = 3280.5208
$ws.Range("K32").Value = 3280.5208
$ws.Range("M32").Value = -2993.5208

$ws.Range("H74").Value = 1675.1923
$ws.Range("I74").Value = 1111.3478
$ws.Range("K74").Value = 1111.3478
$ws.Range("M74").Value = -237.3478

$ws.Range("H77").Value = 1675.1923
$ws.Range("I77").Value = 1111.3478
$ws.Range("K77").Value = 5556.739
$ws.Range("M77").Value = -1188.739

$ws.Range("H80").Value = 59282.332
$ws.Range("I80").Value = 49399
$ws.Range("J80").Value = 61259
$ws.Range("K80").Value = 49399
$ws.Range("L80").Value = 61259
$ws.Range("M80").Value = -48401
$ws.Range("N80").Value = -63255

$ws.Range("H83").Value = 59282.332
$ws.Range("I83").Value = 49399
$ws.Range("J83").Value = 61259
$ws.Range("K83").Value = 148197
$ws.Range("L83").Value = 183777
$ws.Range("M83").Value = -143205
$ws.Range("N83").Value = -193761

$ws.Range("H95").Value = 25850
$ws.Range("J95").Value = 25850
$ws.Range("L95").Value = 25850
$ws.Range("N95").Value = -31342

$ws.Range("H132").Value = 2523.9558
$ws.Range("I132").Value = 2260.5
$ws.Range("K132").Value = 6781.5
$ws.Range("M132").Value = -4251.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 94
$ws.Range("I4").Value = 92.28570999999999
$ws.Range("K4").Value = 92.28570999999999
$ws.Range("M4").Value = 22.71429000000001

$ws.Range("H86").Value = 2432.8518
$ws.Range("I86").Value = 2200.05
$ws.Range("K86").Value = 2200.05
$ws.Range("M86").Value = -1077.05

$ws.Range("H89").Value = 2432.8518
$ws.Range("I89").Value = 2200.05
$ws.Range("K89").Value = 11000.25
$ws.Range("M89").Value = -5384.25

$ws.Range("H99").Value = 1803.12
$ws.Range("I99").Value = 1484.45
$ws.Range("K99").Value = 1484.45
$ws.Range("M99").Value = 13.54999999999995

$ws.Range("H107").Value = 836.86664
$ws.Range("I107").Value = 824.2857
$ws.Range("K107").Value = 824.2857
$ws.Range("M107").Value = 1095.7143

$ws.Range("H137").Value = 66666.664
$ws.Range("J137").Value = 66666.664
$ws.Range("L137").Value = 66666.664
$ws.Range("N137").Value = -76866.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2872.1177
$ws.Range("I31").Value = 1528.25
$ws.Range("K31").Value = 1528.25
$ws.Range("M31").Value = -1233.25

$ws.Range("H34").Value = 2872.1177
$ws.Range("I34").Value = 1528.25
$ws.Range("K34").Value = 1528.25
$ws.Range("M34").Value = -1326.25

$ws.Range("H52").Value = 65135.25
$ws.Range("I52").Value = 53513.668
$ws.Range("J52").Value = 100000
$ws.Range("K52").Value = 53513.668
$ws.Range("L52").Value = 100000
$ws.Range("M52").Value = -53219.668
$ws.Range("N52").Value = -100588

$ws.Range("H62").Value = 5149.5
$ws.Range("J62").Value = 5800
$ws.Range("L62").Value = 5800
$ws.Range("N62").Value = -7048

$ws.Range("H65").Value = 5149.5
$ws.Range("J65").Value = 5800
$ws.Range("L65").Value = 29000
$ws.Range("N65").Value = -35240

$ws.Range("H134").Value = 3769.0513
$ws.Range("I134").Value = 3278.5557
$ws.Range("K134").Value = 9835.667099999999
$ws.Range("M134").Value = -7300.667099999999

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.28571
$ws.Range("I2").Value = 29.2
$ws.Range("K2").Value = 175.2
$ws.Range("M2").Value = -62.19999999999999

$ws.Range("H11").Value = 2222294.5
$ws.Range("I11").Value = 2500075
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 7500225
$ws.Range("L11").Value = 150
$ws.Range("M11").Value = -7500085
$ws.Range("N11").Value = -430

$ws.Range("H101").Value = 11869.571
$ws.Range("J101").Value = 11869.571
$ws.Range("L101").Value = 35608.713
$ws.Range("N101").Value = -40476.713

$ws.Range("H112").Value = 2855.3333
$ws.Range("I112").Value = 2326.6
$ws.Range("K112").Value = 6979.799999999999
$ws.Range("M112").Value = -5871.799999999999

$ws.Range("H117").Value = 943.2857
$ws.Range("I117").Value = 267.16666
$ws.Range("K117").Value = 801.4999799999999
$ws.Range("M117").Value = 2640.50002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H102").Value = 4404.231
$ws.Range("I102").Value = 4727.0293
$ws.Range("J102").Value = 2209.2
$ws.Range("K102").Value = 4727.0293
$ws.Range("L102").Value = 2209.2
$ws.Range("M102").Value = -3105.0293
$ws.Range("N102").Value = -5453.2

$ws.Range("H122").Value = 6702.0938
$ws.Range("I122").Value = 6042.174
$ws.Range("J122").Value = 8388.556
$ws.Range("K122").Value = 18126.522
$ws.Range("L122").Value = 25165.668
$ws.Range("M122").Value = -15676.522
$ws.Range("N122").Value = -30065.668

$ws.Range("H140").Value = 71250
$ws.Range("J140").Value = 71250
$ws.Range("L140").Value = 71250
$ws.Range("N140").Value = -81610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3522.4688
$ws.Range("I46").Value = 3185.5
$ws.Range("J46").Value = 4084.0833
$ws.Range("K46").Value = 3185.5
$ws.Range("L46").Value = 4084.0833
$ws.Range("M46").Value = -2997.5
$ws.Range("N46").Value = -4460.0833

$ws.Range("H61").Value = 3954.52
$ws.Range("I61").Value = 3048.2778
$ws.Range("J61").Value = 6284.857
$ws.Range("K61").Value = 3048.2778
$ws.Range("L61").Value = 6284.857
$ws.Range("M61").Value = -2846.2778
$ws.Range("N61").Value = -6688.857

$ws.Range("H113").Value = 3954.52
$ws.Range("I113").Value = 3048.2778
$ws.Range("J113").Value = 6284.857
$ws.Range("K113").Value = 3048.2778
$ws.Range("L113").Value = 6284.857
$ws.Range("M113").Value = -878.2777999999998
$ws.Range("N113").Value = -10624.857

$ws.Range("H122").Value = 1003385.8
$ws.Range("I122").Value = 1002908.7
$ws.Range("J122").Value = 1003819.56
$ws.Range("K122").Value = 3008726.1
$ws.Range("L122").Value = 3011458.68
$ws.Range("M122").Value = -3006276.1
$ws.Range("N122").Value = -3016358.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11757
$ws.Range("I41").Value = 11578
$ws.Range("K41").Value = 11578
$ws.Range("M41").Value = -11188

$ws.Range("H62").Value = 9240.6
$ws.Range("I62").Value = 8640
$ws.Range("K62").Value = 8640
$ws.Range("M62").Value = -8016

$ws.Range("H65").Value = 9240.6
$ws.Range("I65").Value = 8640
$ws.Range("K65").Value = 43200
$ws.Range("M65").Value = -40080

$ws.Range("H122").Value = 50004420
$ws.Range("J122").Value = 5986.875
$ws.Range("L122").Value = 17960.625
$ws.Range("N122").Value = -22860.625
